$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain plain text (avoids Excel auto-converting
    # numeric-looking strings like "209.11" into floating point numbers),
    # then restore the default "Normal" style so no stray NumberFormat
    # ends up attached to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.157.80"
$ws.Range("E2").Value = "  -2.14%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.578.05"
$ws.Range("E3").Value = "  -1.49%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.33%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "209.11"
$ws.Range("E5").Value = "  -1.29%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.498"
$ws.Range("E6").Value = "  -3.10%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.32%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.246"
$ws.Range("E8").Value = "  -0.67%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.58%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "19.52"
$ws.Range("E10").Value = "  -0.81%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0844"
$ws.Range("E11").Value = "  -0.23%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.799.67"
$ws.Range("E12").Value = "  -1.46%  "

# Row 13 & 14 swap: WrappedEther <-> Polkadot (content swapped rows, with updated price/volume)
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.04"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.565.87"
$ws.Range("E14").Value = "  -2.33%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.514"
$ws.Range("E15").Value = "  -1.90%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "64.44"
$ws.Range("E16").Value = "  -0.98%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.153.33"
$ws.Range("E17").Value = "  -2.04%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -1.09%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "7.24"
$ws.Range("E19").Value = "  +1.20%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "207.85"
$ws.Range("E20").Value = "  -1.03%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.35%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.08%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  -2.91%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  -1.19%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "143.84"
$ws.Range("E25").Value = "  +0.16%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.32%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "6.97"
$ws.Range("E27").Value = "  -1.61%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -1.60%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "15.22"
$ws.Range("E29").Value = "  -0.78%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -0.77%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.53%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "3.21"
$ws.Range("E32").Value = "  -1.62%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "2.97"

# Row 34 - Maker
$ws.Range("D34").Value = "1.279.37"
$ws.Range("E34").Value = "  -0.84%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.50%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +1.25%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -0.98%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -2.53%  "

# Row 39 - WEMIXToken
$ws.Range("E39").Value = "  -5.61%  "

# Row 40 - ARBITRUM
Set-TextValue $ws.Range("D40") "0.817"
$ws.Range("E40").Value = "  -1.67%  "

# Row 41 - FraxShare
$ws.Range("E41").Value = "  +3.06%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  -2.28%  "

# Row 43 & 44 swap: TrustWalletToken <-> Aave
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D43") "62.41"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D44") "0.762"
$ws.Range("E44").Value = "  -2.89%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.712.97"
$ws.Range("E45").Value = "  -1.41%  "

# Row 46 - Quant
Set-TextValue $ws.Range("D46") "88.89"
$ws.Range("E46").Value = "  -1.84%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -0.42%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.23%  "

# Row 49 - Algorand
Set-TextValue $ws.Range("D49") "0.101"
$ws.Range("E49").Value = "  -1.24%  "

# Row 50 - Cronos
Set-TextValue $ws.Range("D50") "0.0505"
$ws.Range("E50").Value = "  -2.10%  "

# Row 51 - USDD -> Aptos (content replaced entirely)
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D51") "5.75"
$ws.Range("E51").Value = "  +11.06%  "
